# ---------------------------------------------------------------------------
# Updates the "cryptos" worksheet with refreshed price / 1h-volume figures
# (and restores the correct coin ordering for the two swapped row pairs),
# matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# Column layout: A=rank(0-based) B=Coin C=Link D=Price E=Volume(1h)
# Price/Volume columns are stored as literal text in the sheet (e.g.
# "72.224.76", "  +0.14%  "), so values that *look* purely numeric (like
# "597.33" or "0.999") need to be written while the cell is temporarily
# formatted as Text ("@") - otherwise Excel would silently reinterpret them
# as numbers and drop formatting/leading characters. The NumberFormat is
# reset back to "Normal" immediately afterwards so no stray cell styling
# is introduced.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [pscustomobject]@{ Row = 2; Col = 4; Value = "72.224.76"; ForceText = $false },
    [pscustomobject]@{ Row = 2; Col = 5; Value = "  +0.14%  "; ForceText = $false },
    [pscustomobject]@{ Row = 3; Col = 4; Value = "2.661.54"; ForceText = $false },
    [pscustomobject]@{ Row = 3; Col = 5; Value = "  +1.71%  "; ForceText = $false },
    [pscustomobject]@{ Row = 4; Col = 5; Value = "  +0.03%  "; ForceText = $false },
    [pscustomobject]@{ Row = 5; Col = 4; Value = "597.33"; ForceText = $true },
    [pscustomobject]@{ Row = 5; Col = 5; Value = "  -0.92%  "; ForceText = $false },
    [pscustomobject]@{ Row = 6; Col = 4; Value = "175.13"; ForceText = $true },
    [pscustomobject]@{ Row = 6; Col = 5; Value = "  -1.51%  "; ForceText = $false },
    [pscustomobject]@{ Row = 8; Col = 5; Value = "  +0.00%  "; ForceText = $false },
    [pscustomobject]@{ Row = 9; Col = 4; Value = "2.660.51"; ForceText = $false },
    [pscustomobject]@{ Row = 9; Col = 5; Value = "  +1.73%  "; ForceText = $false },
    [pscustomobject]@{ Row = 10; Col = 5; Value = "  -1.55%  "; ForceText = $false },
    [pscustomobject]@{ Row = 11; Col = 5; Value = "  +2.49%  "; ForceText = $false },
    [pscustomobject]@{ Row = 12; Col = 5; Value = "  +1.40%  "; ForceText = $false },
    [pscustomobject]@{ Row = 13; Col = 5; Value = "  -0.40%  "; ForceText = $false },
    [pscustomobject]@{ Row = 14; Col = 4; Value = "3.149.23"; ForceText = $false },
    [pscustomobject]@{ Row = 14; Col = 5; Value = "  +1.70%  "; ForceText = $false },
    [pscustomobject]@{ Row = 15; Col = 5; Value = "  -1.03%  "; ForceText = $false },
    [pscustomobject]@{ Row = 16; Col = 4; Value = "71.995.98"; ForceText = $false },
    [pscustomobject]@{ Row = 16; Col = 5; Value = "  -0.06%  "; ForceText = $false },
    [pscustomobject]@{ Row = 17; Col = 4; Value = "26.25"; ForceText = $true },
    [pscustomobject]@{ Row = 17; Col = 5; Value = "  -0.75%  "; ForceText = $false },
    [pscustomobject]@{ Row = 18; Col = 4; Value = "2.654.95"; ForceText = $false },
    [pscustomobject]@{ Row = 18; Col = 5; Value = "  +1.77%  "; ForceText = $false },
    [pscustomobject]@{ Row = 19; Col = 4; Value = "12.24"; ForceText = $true },
    [pscustomobject]@{ Row = 19; Col = 5; Value = "  +6.06%  "; ForceText = $false },
    [pscustomobject]@{ Row = 20; Col = 4; Value = "8.24"; ForceText = $true },
    [pscustomobject]@{ Row = 20; Col = 5; Value = "  +4.61%  "; ForceText = $false },
    [pscustomobject]@{ Row = 21; Col = 4; Value = "370.20"; ForceText = $true },
    [pscustomobject]@{ Row = 21; Col = 5; Value = "  -3.10%  "; ForceText = $false },
    [pscustomobject]@{ Row = 22; Col = 5; Value = "  +0.30%  "; ForceText = $false },
    [pscustomobject]@{ Row = 23; Col = 5; Value = "  +2.64%  "; ForceText = $false },
    [pscustomobject]@{ Row = 24; Col = 4; Value = "71.98"; ForceText = $true },
    [pscustomobject]@{ Row = 24; Col = 5; Value = "  -1.32%  "; ForceText = $false },
    [pscustomobject]@{ Row = 25; Col = 5; Value = "  +0.18%  "; ForceText = $false },
    [pscustomobject]@{ Row = 26; Col = 4; Value = "4.32"; ForceText = $true },
    [pscustomobject]@{ Row = 26; Col = 5; Value = "  -0.91%  "; ForceText = $false },
    [pscustomobject]@{ Row = 27; Col = 4; Value = "9.76"; ForceText = $true },
    [pscustomobject]@{ Row = 27; Col = 5; Value = "  -0.68%  "; ForceText = $false },
    [pscustomobject]@{ Row = 28; Col = 4; Value = "2.797.67"; ForceText = $false },
    [pscustomobject]@{ Row = 28; Col = 5; Value = "  +1.70%  "; ForceText = $false },
    [pscustomobject]@{ Row = 29; Col = 4; Value = "0.999"; ForceText = $true },
    [pscustomobject]@{ Row = 29; Col = 5; Value = "  -0.11%  "; ForceText = $false },
    [pscustomobject]@{ Row = 30; Col = 4; Value = "0.0₃0968"; ForceText = $false },
    [pscustomobject]@{ Row = 30; Col = 5; Value = "  +2.33%  "; ForceText = $false },
    [pscustomobject]@{ Row = 31; Col = 4; Value = "8.06"; ForceText = $true },
    [pscustomobject]@{ Row = 31; Col = 5; Value = "  +0.63%  "; ForceText = $false },
    [pscustomobject]@{ Row = 32; Col = 4; Value = "500.92"; ForceText = $true },
    [pscustomobject]@{ Row = 32; Col = 5; Value = "  -3.37%  "; ForceText = $false },
    [pscustomobject]@{ Row = 33; Col = 5; Value = "  -1.94%  "; ForceText = $false },
    [pscustomobject]@{ Row = 34; Col = 5; Value = "  +0.13%  "; ForceText = $false },
    [pscustomobject]@{ Row = 35; Col = 5; Value = "  +0.02%  "; ForceText = $false },
    [pscustomobject]@{ Row = 36; Col = 4; Value = "162.99"; ForceText = $true },
    [pscustomobject]@{ Row = 36; Col = 5; Value = "  -0.34%  "; ForceText = $false },
    [pscustomobject]@{ Row = 37; Col = 4; Value = "19.51"; ForceText = $true },
    [pscustomobject]@{ Row = 37; Col = 5; Value = "  +1.43%  "; ForceText = $false },
    [pscustomobject]@{ Row = 38; Col = 5; Value = "  +0.05%  "; ForceText = $false },
    [pscustomobject]@{ Row = 39; Col = 4; Value = "18.97"; ForceText = $true },
    [pscustomobject]@{ Row = 39; Col = 5; Value = "  -0.57%  "; ForceText = $false },
    [pscustomobject]@{ Row = 40; Col = 5; Value = "  -1.35%  "; ForceText = $false },
    [pscustomobject]@{ Row = 41; Col = 5; Value = "  -2.75%  "; ForceText = $false },
    [pscustomobject]@{ Row = 42; Col = 5; Value = "  +0.06%  "; ForceText = $false },
    [pscustomobject]@{ Row = 43; Col = 5; Value = "  -0.83%  "; ForceText = $false },
    [pscustomobject]@{ Row = 44; Col = 2; Value = "PolygonEcosystemToken"; ForceText = $false },
    [pscustomobject]@{ Row = 44; Col = 3; Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; ForceText = $false },
    [pscustomobject]@{ Row = 44; Col = 4; Value = "0.332"; ForceText = $true },
    [pscustomobject]@{ Row = 44; Col = 5; Value = "  +0.48%  "; ForceText = $false },
    [pscustomobject]@{ Row = 45; Col = 2; Value = "dogwifhat"; ForceText = $false },
    [pscustomobject]@{ Row = 45; Col = 3; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; ForceText = $false },
    [pscustomobject]@{ Row = 45; Col = 4; Value = "2.55"; ForceText = $true },
    [pscustomobject]@{ Row = 45; Col = 5; Value = "  -0.76%  "; ForceText = $false },
    [pscustomobject]@{ Row = 46; Col = 4; Value = "156.61"; ForceText = $true },
    [pscustomobject]@{ Row = 46; Col = 5; Value = "  +4.36%  "; ForceText = $false },
    [pscustomobject]@{ Row = 47; Col = 4; Value = "39.47"; ForceText = $true },
    [pscustomobject]@{ Row = 47; Col = 5; Value = "  +0.00%  "; ForceText = $false },
    [pscustomobject]@{ Row = 48; Col = 2; Value = "Filecoin"; ForceText = $false },
    [pscustomobject]@{ Row = 48; Col = 3; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; ForceText = $false },
    [pscustomobject]@{ Row = 48; Col = 4; Value = "3.74"; ForceText = $true },
    [pscustomobject]@{ Row = 48; Col = 5; Value = "  +1.76%  "; ForceText = $false },
    [pscustomobject]@{ Row = 49; Col = 2; Value = "ARBITRUM"; ForceText = $false },
    [pscustomobject]@{ Row = 49; Col = 3; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; ForceText = $false },
    [pscustomobject]@{ Row = 49; Col = 4; Value = "0.558"; ForceText = $true },
    [pscustomobject]@{ Row = 49; Col = 5; Value = "  +3.24%  "; ForceText = $false },
    [pscustomobject]@{ Row = 50; Col = 5; Value = "  +2.20%  "; ForceText = $false },
    [pscustomobject]@{ Row = 51; Col = 4; Value = "0.0754"; ForceText = $true },
    [pscustomobject]@{ Row = 51; Col = 5; Value = "  -1.42%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Output ("Applied {0} cell updates to Sheet1." -f $updates.Count)
